$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text numeric-looking values (e.g. "70.00", "41.618.42").
# Force text format first so Excel does not auto-convert them to numbers,
# which would corrupt values like multi-dot numbers or drop trailing zeros.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '41.618.42'
$ws.Range('E2').Value = '  -1.49%  '
$ws.Range('D3').Value = '2.165.06'
$ws.Range('E3').Value = '  -2.95%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '237.87'
$ws.Range('D6').Value = '0.609'
$ws.Range('E6').Value = '  -2.87%  '
$ws.Range('D7').Value = '72.22'
$ws.Range('E7').Value = '  -2.83%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.579'
$ws.Range('E9').Value = '  -4.34%  '
$ws.Range('D10').Value = '39.38'
$ws.Range('E10').Value = '  -7.90%  '
$ws.Range('E11').Value = '  -5.66%  '
$ws.Range('D12').Value = '54.38'
$ws.Range('E12').Value = '  -3.91%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '6.69'
$ws.Range('E13').Value = '  -3.96%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').Value = '0.0997'
$ws.Range('E14').Value = '  -3.73%  '
$ws.Range('D15').Value = '2.490.20'
$ws.Range('E15').Value = '  -2.91%  '
$ws.Range('D16').Value = '14.28'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = '2.148.69'
$ws.Range('E17').Value = '  -4.08%  '
$ws.Range('D18').Value = '0.777'
$ws.Range('E18').Value = '  -7.31%  '
$ws.Range('D19').Value = '41.506.14'
$ws.Range('E19').Value = '  -1.31%  '
$ws.Range('D20').Value = '0.0000102'
$ws.Range('E20').Value = '  -3.26%  '
$ws.Range('D21').Value = '70.00'
$ws.Range('E21').Value = '  -3.93%  '
$ws.Range('D22').Value = '5.78'
$ws.Range('E22').Value = '  -7.14%  '
$ws.Range('D23').Value = '9.91'
$ws.Range('E23').Value = '  -11.12%  '
$ws.Range('D24').Value = '226.21'
$ws.Range('E24').Value = '  -1.96%  '
$ws.Range('E25').Value = '  -4.20%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = '10.71'
$ws.Range('E27').Value = '  -5.89%  '
$ws.Range('D29').Value = '2.18'
$ws.Range('E29').Value = '  -4.06%  '
$ws.Range('E30').Value = '  -1.71%  '
$ws.Range('D31').Value = '171.10'
$ws.Range('E31').Value = '  +2.46%  '
$ws.Range('D32').Value = '19.81'
$ws.Range('E32').Value = '  -3.95%  '
$ws.Range('D33').Value = '32.77'
$ws.Range('E33').Value = '  +9.04%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '5.40'
$ws.Range('E34').Value = '  -4.44%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.0771'
$ws.Range('E35').Value = '  -4.16%  '
$ws.Range('D36').Value = '0.120'
$ws.Range('E36').Value = '  -3.71%  '
$ws.Range('E37').Value = '  -1.13%  '
$ws.Range('D38').Value = '0.103'
$ws.Range('E38').Value = '  -7.00%  '
$ws.Range('E39').Value = '  -0.28%  '
$ws.Range('D40').Value = '12.17'
$ws.Range('E40').Value = '  -7.92%  '
$ws.Range('E41').Value = '  -1.95%  '
$ws.Range('E42').Value = '  -6.05%  '
$ws.Range('D43').Value = '58.78'
$ws.Range('E43').Value = '  -9.80%  '
$ws.Range('D44').Value = '8.45'
$ws.Range('E44').Value = '  -3.18%  '
$ws.Range('E45').Value = '  -5.43%  '
$ws.Range('D46').Value = '0.0963'
$ws.Range('E46').Value = '  -3.93%  '
$ws.Range('D47').Value = '97.36'
$ws.Range('E47').Value = '  -6.96%  '
$ws.Range('E48').Value = '  -4.02%  '
$ws.Range('E49').Value = '  -5.01%  '
$ws.Range('E50').Value = '  -7.41%  '
$ws.Range('E51').Value = '  -2.38%  '

# Restore original (default) formatting on column D now that the text values are set.
$dRange.ClearFormats()

